# The deck has three "title only" slides (VISION, SOFTWARE ARCHITECTURE,
# CHALLENGES & RISKS) that each immediately precede a matching slide that
# repeats the same title together with the actual bullet content. This
# edit removes the redundant title-only slides, keeping only the slides
# that carry the real content (their titles already read the same).
#
# Original slide order (1-indexed):
#   1  PROJECT:CORNUCOPIA
#   2  OUR TEAM
#   3  TABLE OF CONTENTS
#   4  INTRODUCTION
#   5  VISION                      (title only)               <- remove
#   6  VISION                      (title + bullets)           <- keep
#   7  SOFTWARE ARCHITECTURE       (title only)               <- remove
#   8  SOFTWARE ARCHITECTURE       (title + bullets)           <- keep
#   9  CHALLENGES & RISKS          (title only)               <- remove
#   10 CHALLENGES & RISKS          (title + bullets)           <- keep
#
# Deleting from the highest index down avoids any index shifting issues.

$p = $ppt.ActivePresentation

$p.Slides.Item(9).Delete()
$p.Slides.Item(7).Delete()
$p.Slides.Item(5).Delete()
